# Scheduled-runner style refresh of cached market-board stats (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Leve sheets. No formulas are involved -
# these columns are plain cached numeric snapshots - so we just overwrite the literal
# values per affected row on each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Cells.Item(41, 8).Value = 219.53847
$ws.Cells.Item(41, 9).Value = 241.27272
$ws.Cells.Item(41, 10).Value = 100
$ws.Cells.Item(41, 11).Value = 241.27272
$ws.Cells.Item(41, 12).Value = 100
$ws.Cells.Item(41, 13).Value = 198.72728
$ws.Cells.Item(41, 14).Value = -980

# Row 100
$ws.Cells.Item(100, 8).Value = 1343.0769
$ws.Cells.Item(100, 9).Value = 1000
$ws.Cells.Item(100, 11).Value = 1000
$ws.Cells.Item(100, 13).Value = -459

# Row 113
$ws.Cells.Item(113, 8).Value = 1711.4286
$ws.Cells.Item(113, 10).Value = 1711.4286
$ws.Cells.Item(113, 12).Value = 1711.4286
$ws.Cells.Item(113, 14).Value = -8219.428599999999

# Row 116
$ws.Cells.Item(116, 8).Value = 2939.4
$ws.Cells.Item(116, 9).Value = 3627.7144
$ws.Cells.Item(116, 10).Value = 1333.3334
$ws.Cells.Item(116, 11).Value = 3627.7144
$ws.Cells.Item(116, 12).Value = 1333.3334
$ws.Cells.Item(116, 13).Value = -185.7143999999998
$ws.Cells.Item(116, 14).Value = -8217.3334

# Row 129
$ws.Cells.Item(129, 8).Value = 1076.75
$ws.Cells.Item(129, 9).Value = 461.53845
$ws.Cells.Item(129, 10).Value = 1803.8182
$ws.Cells.Item(129, 11).Value = 1384.61535
$ws.Cells.Item(129, 12).Value = 5411.4546
$ws.Cells.Item(129, 13).Value = 3615.38465
$ws.Cells.Item(129, 14).Value = -15411.4546

# Row 132
$ws.Cells.Item(132, 8).Value = 259421.28
$ws.Cells.Item(132, 9).Value = 297404.94
$ws.Cells.Item(132, 10).Value = 1132.4
$ws.Cells.Item(132, 11).Value = 892214.8200000001
$ws.Cells.Item(132, 12).Value = 3397.2
$ws.Cells.Item(132, 13).Value = -889684.8200000001
$ws.Cells.Item(132, 14).Value = -8457.200000000001

# Row 138
$ws.Cells.Item(138, 8).Value = 2293.2034
$ws.Cells.Item(138, 9).Value = 1883.7941
$ws.Cells.Item(138, 10).Value = 2850
$ws.Cells.Item(138, 11).Value = 5651.3823
$ws.Cells.Item(138, 12).Value = 8550
$ws.Cells.Item(138, 13).Value = -511.3823000000002
$ws.Cells.Item(138, 14).Value = -18830

# Row 141
$ws.Cells.Item(141, 8).Value = 990.1142599999999
$ws.Cells.Item(141, 9).Value = 763.931
$ws.Cells.Item(141, 10).Value = 2083.3333
$ws.Cells.Item(141, 11).Value = 2291.793
$ws.Cells.Item(141, 12).Value = 6249.999899999999
$ws.Cells.Item(141, 13).Value = 2888.207
$ws.Cells.Item(141, 14).Value = -16609.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 3895.6628
$ws.Cells.Item(32, 9).Value = 3899.4692
$ws.Cells.Item(32, 10).Value = 3834
$ws.Cells.Item(32, 11).Value = 3899.4692
$ws.Cells.Item(32, 12).Value = 3834
$ws.Cells.Item(32, 13).Value = -3612.4692
$ws.Cells.Item(32, 14).Value = -4408

# Row 61
$ws.Cells.Item(61, 8).Value = 9526285
$ws.Cells.Item(61, 9).Value = 10419123
$ws.Cells.Item(61, 10).Value = 2679.3333
$ws.Cells.Item(61, 11).Value = 10419123
$ws.Cells.Item(61, 12).Value = 2679.3333
$ws.Cells.Item(61, 13).Value = -10418911
$ws.Cells.Item(61, 14).Value = -3103.3333

# Row 74
$ws.Cells.Item(74, 8).Value = 5409975
$ws.Cells.Item(74, 9).Value = 7408070
$ws.Cells.Item(74, 10).Value = 15117.7
$ws.Cells.Item(74, 11).Value = 7408070
$ws.Cells.Item(74, 12).Value = 15117.7
$ws.Cells.Item(74, 13).Value = -7407196
$ws.Cells.Item(74, 14).Value = -16865.7

# Row 77
$ws.Cells.Item(77, 8).Value = 5409975
$ws.Cells.Item(77, 9).Value = 7408070
$ws.Cells.Item(77, 10).Value = 15117.7
$ws.Cells.Item(77, 11).Value = 37040350
$ws.Cells.Item(77, 12).Value = 75588.5
$ws.Cells.Item(77, 13).Value = -37035982
$ws.Cells.Item(77, 14).Value = -84324.5

# Row 122
$ws.Cells.Item(122, 8).Value = 1514.5161
$ws.Cells.Item(122, 9).Value = 1496.5927
$ws.Cells.Item(122, 11).Value = 4489.7781
$ws.Cells.Item(122, 13).Value = -2039.7781

# Row 132
$ws.Cells.Item(132, 8).Value = 711139.1
$ws.Cells.Item(132, 9).Value = 735186
$ws.Cells.Item(132, 10).Value = 334404.66
$ws.Cells.Item(132, 11).Value = 2205558
$ws.Cells.Item(132, 12).Value = 1003213.98
$ws.Cells.Item(132, 13).Value = -2203028
$ws.Cells.Item(132, 14).Value = -1008273.98

# Row 136
$ws.Cells.Item(136, 8).Value = 9526285
$ws.Cells.Item(136, 9).Value = 10419123
$ws.Cells.Item(136, 10).Value = 2679.3333
$ws.Cells.Item(136, 11).Value = 31257369
$ws.Cells.Item(136, 12).Value = 8037.999899999999
$ws.Cells.Item(136, 13).Value = -31254819
$ws.Cells.Item(136, 14).Value = -13137.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Cells.Item(59, 8).Value = 33995
$ws.Cells.Item(59, 10).Value = 33995
$ws.Cells.Item(59, 12).Value = 33995
$ws.Cells.Item(59, 14).Value = -35689

# Row 94
$ws.Cells.Item(94, 8).Value = 690.8182
$ws.Cells.Item(94, 9).Value = 514.1429000000001
$ws.Cells.Item(94, 10).Value = 1000
$ws.Cells.Item(94, 11).Value = 514.1429000000001
$ws.Cells.Item(94, 12).Value = 1000
$ws.Cells.Item(94, 13).Value = -63.14290000000005
$ws.Cells.Item(94, 14).Value = -1902

# Row 99
$ws.Cells.Item(99, 8).Value = 1228.5
$ws.Cells.Item(99, 9).Value = 1020.75
$ws.Cells.Item(99, 11).Value = 1020.75
$ws.Cells.Item(99, 13).Value = 477.25

# Row 134
$ws.Cells.Item(134, 8).Value = 81079.336
$ws.Cells.Item(134, 9).Value = 86726.86
$ws.Cells.Item(134, 10).Value = 2014
$ws.Cells.Item(134, 11).Value = 260180.58
$ws.Cells.Item(134, 12).Value = 6042
$ws.Cells.Item(134, 13).Value = -257645.58
$ws.Cells.Item(134, 14).Value = -11112

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Cells.Item(99, 8).Value = 2011.1578
$ws.Cells.Item(99, 9).Value = 2114
$ws.Cells.Item(99, 10).Value = 1936.3636
$ws.Cells.Item(99, 11).Value = 2114
$ws.Cells.Item(99, 12).Value = 1936.3636
$ws.Cells.Item(99, 13).Value = -616
$ws.Cells.Item(99, 14).Value = -4932.3636

# Row 109
$ws.Cells.Item(109, 8).Value = 40000
$ws.Cells.Item(109, 10).Value = 40000
$ws.Cells.Item(109, 12).Value = 40000
$ws.Cells.Item(109, 14).Value = -42080

# Row 126
$ws.Cells.Item(126, 8).Value = 2011.1578
$ws.Cells.Item(126, 9).Value = 2114
$ws.Cells.Item(126, 10).Value = 1936.3636
$ws.Cells.Item(126, 11).Value = 6342
$ws.Cells.Item(126, 12).Value = 5809.0908
$ws.Cells.Item(126, 13).Value = -3872
$ws.Cells.Item(126, 14).Value = -10749.0908

# Row 132
$ws.Cells.Item(132, 8).Value = 1643.8541
$ws.Cells.Item(132, 9).Value = 1703.375
$ws.Cells.Item(132, 10).Value = 1346.25
$ws.Cells.Item(132, 11).Value = 5110.125
$ws.Cells.Item(132, 12).Value = 4038.75
$ws.Cells.Item(132, 13).Value = -2580.125
$ws.Cells.Item(132, 14).Value = -9098.75

# Row 134
$ws.Cells.Item(134, 8).Value = 2857.2144
$ws.Cells.Item(134, 9).Value = 2991.3076
$ws.Cells.Item(134, 10).Value = 1114
$ws.Cells.Item(134, 11).Value = 8973.9228
$ws.Cells.Item(134, 12).Value = 3342
$ws.Cells.Item(134, 13).Value = -6438.9228
$ws.Cells.Item(134, 14).Value = -8412

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Cells.Item(32, 8).Value = 24290
$ws.Cells.Item(32, 10).Value = 24290
$ws.Cells.Item(32, 12).Value = 24290
$ws.Cells.Item(32, 14).Value = -24882

# Row 42
$ws.Cells.Item(42, 8).Value = 28516
$ws.Cells.Item(42, 10).Value = 28516
$ws.Cells.Item(42, 12).Value = 28516
$ws.Cells.Item(42, 14).Value = -29486

# Row 115
$ws.Cells.Item(115, 8).Value = 28516
$ws.Cells.Item(115, 10).Value = 28516
$ws.Cells.Item(115, 12).Value = 28516
$ws.Cells.Item(115, 14).Value = -30866

# Row 132
$ws.Cells.Item(132, 8).Value = 1576.1305
$ws.Cells.Item(132, 9).Value = 1881.125
$ws.Cells.Item(132, 10).Value = 1243.409
$ws.Cells.Item(132, 11).Value = 5643.375
$ws.Cells.Item(132, 12).Value = 3730.227
$ws.Cells.Item(132, 13).Value = -3113.375
$ws.Cells.Item(132, 14).Value = -8790.227000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Cells.Item(136, 8).Value = 4530.0645
$ws.Cells.Item(136, 9).Value = 4083.862
$ws.Cells.Item(136, 10).Value = 11000
$ws.Cells.Item(136, 11).Value = 12251.586
$ws.Cells.Item(136, 12).Value = 33000
$ws.Cells.Item(136, 13).Value = -9701.585999999999
$ws.Cells.Item(136, 14).Value = -38100

# Row 138
$ws.Cells.Item(138, 8).Value = 64645.8
$ws.Cells.Item(138, 10).Value = 64645.8
$ws.Cells.Item(138, 12).Value = 64645.8
$ws.Cells.Item(138, 14).Value = -74925.8

$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Cells.Item(109, 8).Value = 49377
$ws.Cells.Item(109, 10).Value = 49377
$ws.Cells.Item(109, 12).Value = 49377
$ws.Cells.Item(109, 14).Value = -52151

# Row 132
$ws.Cells.Item(132, 8).Value = 5340.394
$ws.Cells.Item(132, 9).Value = 6224.4287
$ws.Cells.Item(132, 10).Value = 389.8
$ws.Cells.Item(132, 11).Value = 18673.2861
$ws.Cells.Item(132, 12).Value = 1169.4
$ws.Cells.Item(132, 13).Value = -16143.2861
$ws.Cells.Item(132, 14).Value = -6229.4
